$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking scrape refresh: update Price (D) and 1h Volume change (E) columns.
# Pre-format the data range as Text so numeric-looking strings (e.g. "0.590", "1.00",
# "52.010.82") are stored verbatim instead of being re-parsed/normalized as numbers,
# then restore the Normal style so no residual formatting is left on the cells.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '52.010.82'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '2.777.57'
$ws.Range("E3").Value = '  -1.77%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '357.53'
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("D6").Value = '109.36'
$ws.Range("E6").Value = '  -3.99%  '
$ws.Range("D7").Value = '0.563'
$ws.Range("E7").Value = '  +2.79%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").Value = '0.590'
$ws.Range("E9").Value = '  -2.53%  '
$ws.Range("D10").Value = '40.14'
$ws.Range("E10").Value = '  -4.35%  '
$ws.Range("D11").Value = '0.0849'
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("E12").Value = '  +0.74%  '
$ws.Range("D13").Value = '19.43'
$ws.Range("E13").Value = '  -3.50%  '
$ws.Range("E14").Value = '  -2.39%  '
$ws.Range("D15").Value = '3.218.23'
$ws.Range("E15").Value = '  -1.20%  '
$ws.Range("D16").Value = '2.762.34'
$ws.Range("E16").Value = '  -2.16%  '
$ws.Range("D17").Value = '0.921'
$ws.Range("E17").Value = '  +3.38%  '
$ws.Range("D18").Value = '51.900.96'
$ws.Range("E18").Value = '  -0.44%  '
$ws.Range("E19").Value = '  +0.95%  '
$ws.Range("E20").Value = '  -1.09%  '
$ws.Range("D21").Value = '13.08'
$ws.Range("E21").Value = '  -5.38%  '
$ws.Range("D22").Value = '0.0₃0975'
$ws.Range("E22").Value = '  -2.19%  '
$ws.Range("D23").Value = '274.06'
$ws.Range("E23").Value = '  +1.23%  '
$ws.Range("D24").Value = '69.65'
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("E25").Value = '  -1.63%  '
$ws.Range("D26").Value = '26.55'
$ws.Range("E26").Value = '  -0.60%  '
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").Value = '10.12'
$ws.Range("E28").Value = '  -1.42%  '
$ws.Range("D29").Value = '2.24'
$ws.Range("E29").Value = '  -0.74%  '
$ws.Range("E30").Value = '  +2.34%  '
$ws.Range("E31").Value = '  +4.69%  '
$ws.Range("D32").Value = '51.39'
$ws.Range("E32").Value = '  +1.24%  '
$ws.Range("E33").Value = '  +0.64%  '
$ws.Range("D34").Value = '5.71'
$ws.Range("E34").Value = '  -2.99%  '
$ws.Range("D35").Value = '5.41'
$ws.Range("E35").Value = '  +10.86%  '
$ws.Range("D36").Value = '0.0835'
$ws.Range("E36").Value = '  +0.43%  '
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.19%  '
$ws.Range("E38").Value = '  +0.23%  '
$ws.Range("D39").Value = '18.24'
$ws.Range("E39").Value = '  -1.21%  '
$ws.Range("E40").Value = '  -4.57%  '
$ws.Range("D41").Value = '2.54'
$ws.Range("E41").Value = '  -2.16%  '
$ws.Range("E42").Value = '  -0.94%  '
$ws.Range("D43").Value = '123.97'
$ws.Range("E43").Value = '  -2.98%  '
$ws.Range("E44").Value = '  -2.08%  '
$ws.Range("D45").Value = '21.93'
$ws.Range("E45").Value = '  -6.35%  '
$ws.Range("D46").Value = '2.067.11'
$ws.Range("E46").Value = '  +1.20%  '
$ws.Range("E47").Value = '  -3.53%  '
$ws.Range("E48").Value = '  -0.96%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("D50").Value = '0.924'
$ws.Range("E50").Value = '  -4.92%  '
$ws.Range("D51").Value = '8.95'
$ws.Range("E51").Value = '  +0.60%  '

$dataRange.Style = "Normal"
